$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# Remove the empty per-year placeholder sheets; all publication rows
# live in the single consolidated "publications" sheet now.
$sheetsToRemove = @("Sheet9", "vem2017", "vem2018", "vem2019", "vem2020", "vem2021", "vem2022")
foreach ($name in $sheetsToRemove) {
    $wb.Worksheets.Item($name).Delete()
}

# Rename the main sheet to match the consolidated workbook name.
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "publications"

# New VEM 2017 / VEM 2018 publications appended after the existing rows.
$newRows = @(
    @(154, 'Caracterização do Papel Desempenhado por Desenvolvedores Responsáveis pelo Truck Factor de Projetos de Software', 'Thaís Mombach; Mívian Ferreira; Marco Tulio Valente; Kecia Ferreira ', 'UFMG; CEFET-MG', 'vem', 2017),
    @(155, 'Como Pesquisadores Usam o Dataset GHTorrent?', 'Hudson Borges; Jailton Coelho; Paulo Carvalho; Mariane Fernandes; Marco Tulio Valente ', 'UFMG', 'vem', 2017),
    @(156, 'Uma Análise Preliminar da Hiperespecialização no Contexto de Crowdsourcing na Plataforma Topcoder', 'Anderson Bergamini de Neira; Igor Scaliante Wiese; Igor Steinmacher ', 'UEM; UTFPR', 'vem', 2017),
    @(157, 'Who Drive Company-Owned OSS Projects: Employees or Volunteers?', 'Luiz Felipe Dias; Jhoylan Santos; Igor Steinmacher; Gustavo Pinto ', 'UTFPR; UFPA', 'vem', 2017),
    @(158, 'Do Coupling Metrics Help Characterize Critical Components in Component-based SPL? An Empirical Study', 'Anderson Uchôa; Eduardo Fernandes; Ana Carla Bibiano; Alessandro Garcia ', 'PUC-Rio', 'vem', 2017),
    @(159, 'Impactos de Práticas de Codificação sobre a Legibilidade: Um Estudo Preliminar', 'Rodrigo Magalhães Dos Santos; Marco Aurélio Gerosa ', 'IPT-SP; NAU', 'vem', 2017),
    @(160, 'VisminerTD: Uma Ferramenta para Identificação Automática e Monitoramento Interativo de Dívida Técnica', 'Thiago Mendes; David Gonçalves; Felipe Gomes; Renato Novais; Rodrigo Spínola; Manoel Mendonça ', 'UFBA; IFBA; UNIFACS', 'vem', 2017),
    @(161, 'Vulnerabilidades em Sistemas de Software Web baseados em Plug-ins: Um Estudo Exploratório do WordPress', 'Oslien Mesa; Marx Viana; Elder Cirilo; Carlos Lucena ', 'PUC-Rio; UFJF; FIVJ', 'vem', 2017),
    @(162, 'Em Busca de uma Abordagem de Conformidade Arquitetural para Arquitetura de Microsserviços', 'Elena A. Araujo; Elder Rodrigues Jr; Arthur F. Pinto; Ricardo Terra ', 'UFLA', 'vem', 2017),
    @(163, 'Arquitetura de Aplicações Spring MVC: Uma Análise Baseada no Acoplamento Lógico', 'Eduardo Pereira de Sousa; Eduardo Martins Guerra; Gustavo Ansaldi Oliva; Mauricio Finavaro Aniche ', 'INPE; Queen''s University; Delft University', 'vem', 2017),
    @(164, 'Uma Caracterização em Larga Escala da Arquitetura de Sistemas Docker', 'Lucas Monteiro; Laerte Xavier; Marco Tulio Valente ', 'UFMG', 'vem', 2017),
    @(165, 'Usando aprendizagem de máquina para identificar anomalias de design prejudiciais à manutenibilidade: um estudo preliminar', 'Fabricio F Cardim; Cláudio Sant’Anna ', 'UFBA', 'vem', 2017),
    @(166, 'Mineração de regras para solução de problemas relacionados à fragmentação do Android', 'Marcus Adriano F Pereira; Marcelo Maia ', 'UFU', 'vem', 2017),
    @(167, 'Definição de clusters para classificação do uso de anotações em código Java', 'Phyllipe Lima; Eduardo Guerra; Paulo Meirelles ', 'INPE; USP', 'vem', 2017),
    @(168, 'Caracterizando o Consumo de Energia de APIs de E/S da Linguagem Java', 'Gilson Rocha; Gustavo Pinto; Fernando Castor', 'UFPA; UFPE', 'vem', 2018),
    @(169, 'Um Estudo em Larga-Escala sobre Características de APIs Populares', 'Caroline Lima; Pedro Henrique de Moraes; Andre Hora', 'UFMS', 'vem', 2018),
    @(170, 'Minerando Mensagens de Depreciação Faltantes em APIs: Um Estudo de Caso no Ecossistema Android', 'Pedro Henrique de Moraes; Caroline Lima; Andre Hora', 'UFMS', 'vem', 2018),
    @(171, 'STF: uma abordagem Social para estimar Truck Factor no GitHub', 'Hercules Sandim; Michele A. Brandão; Mirella M. Moro', 'UFMG', 'vem', 2018),
    @(172, 'Heurísticas para Identificação de Ambiguidade de Autores em Projetos Open Source', 'Talita Santana Orfano; Kecia Aline Marques Ferreira; Mariza Andrade Da Silva Bigonha', 'UFMG; CEFET-MG', 'vem', 2018),
    @(173, 'Monorepos: A Multivocal Literature Review', 'Gleison Brito; Ricardo Terra; Marco Tulio Valente', 'UFMG; UFLA', 'vem', 2018),
    @(174, 'Minerando Código Comentado', 'Lucas Grijó; Andre Hora', 'UFMS', 'vem', 2018),
    @(175, 'Um Estudo Empírico sobre Critérios de Seleção de Repositórios GitHub', 'Laerte Xavier; Jailton Coelho; Luciana L. Silva', 'UFMG; IFMG', 'vem', 2018),
    @(176, 'GitHub REST API vs GHTorrent vs GitHub Archive: A Comparative Study', 'Thaís Mombach; Marco Tulio Valente', 'UFMG', 'vem', 2018),
    @(177, 'Microservices in Practice: A Survey Study', 'Markos Viggiato; Ricardo Terra; Henrique Rocha; Marco Tulio Valente; Eduardo Figueiredo', 'UFMG; UFLA; Inria', 'vem', 2018),
    @(178, 'Um Método para Detectar Similaridade entre Sistemas baseado em Decisões de Design: um Estudo Preliminar', 'Marcos Dósea; Claudio Sant''Anna', 'UFBA; UFS', 'vem', 2018),
    @(179, 'Comparando Técnicas de Extração de Valores Limiares para Métricas: Um Estudo Preliminar com Desenvolvedores Web', 'Raphael Lima; Marcos Dósea; Cláudio Sant''Anna', 'UFBA; UFS', 'vem', 2018),
    @(180, 'Um Estudo Empírico sobre o Impacto dos Pré-processamentos e Normalizações no Cálculo do Acoplamento Conceitual', 'Paulo Batista Da Costa; Igor Scaliante Wiese; Reginaldo Re; Igor Steinmacher', 'UTFPR', 'vem', 2018),
    @(181, 'Violação de padrões de uso de APIs em sistemas configuráveis', 'Bruno Mecca; Diogo Boaventura; Bruno Cafeo; Elder Cirilo', 'UFMS; UFSJ', 'vem', 2018),
    @(182, 'Avaliação da Frequência de Mudanças em Dependências entre Variabilidades em Sistemas Configuráveis', 'Raiza Artemam de Oliveira; Bruno Gonçalves Mecca; Bruno Barbieri de Pontes Cafeo; Andre Hora', 'UFMS', 'vem', 2018),
    @(183, 'VMAG 3D – An approach for supporting the comprehension of software system models using motion control in a multiuser 3D visualization environment', 'Sergio Henriques Martins Barreto Bento Antunes; Claudia Susie Camargo Rodrigues; Cláudia Maria Lima Werner', 'COPPE/UFRJ', 'vem', 2018),
    @(184, 'Development and Maintenance of Model-Oriented Software with Visualization - Exploratory and Experimental Study', 'Thiago Gottardi; Rosana Braga', 'USP', 'vem', 2018),
    @(185, 'Uma Análise da Produção Científica Brasileira em Conferências de Manutenção e Evolução de Software', 'Klérisson Paixão; Marcelo Maia; Marco Tulio Valente', 'UFU; UFMG', 'vem', 2018),
    @(186, 'An Infrastructure for Software Release Analysis through Provenance Graphs', 'Felipe Curty; Troy Kohwalter; Vanessa Braganholo; Leonardo Murta', 'UFF', 'vem', 2018),
    @(187, 'Uma técnica para a quantificação do esforço de merge', 'Tayane Moura; Leonardo Murta', 'UFF', 'vem', 2018),
    @(188, 'Towards an automated approach for bug fix pattern detection', 'Fernanda Madeiral; Thomas Durieux; Victor Sobreira; Marcelo Maia', 'UFU; Inria', 'vem', 2018),
    @(189, 'Explorando Como Bibliotecas Python Lançam Exceções ao Longo de sua Evolução', 'Allan Gonçalves; Cinthia Nascimento; Eiji Adachi', 'UFRN', 'vem', 2018),
    @(190, 'Identifying Confusing Code in Swift Programs', 'Fernando Castor', 'UFPE', 'vem', 2018),
    @(191, 'DiffMutAnalyze: Uma abordagem para auxiliar a identificação de mutantes equivalentes', 'Juliana Botelho; Carlos Henrique Pereira; Vinicius H. S. Durelli; Rafael S. Durelli', 'UFLA; UFSJ', 'vem', 2018)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}

$ws.Activate()
$ws.Range("G182").Select()
